$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.930.23'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '1.624.04'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '213.83'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.503'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.250'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.23%  '
$ws.Range('E9').Value = '  -3.43%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.28'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -6.44%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.762.69'
$ws.Range('E12').Value = '  +6.93%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.847.81'
$ws.Range('E13').Value = '  -1.14%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.18'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.22%  '
$ws.Range('E15').Value = '  -3.57%  '
$ws.Range('D16').Value = '25.927.78'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.0₃0735'
$ws.Range('E17').Value = '  -3.77%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '61.11'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -3.60%  '
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '191.50'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.43%  '
$ws.Range('E21').Value = '  -3.20%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.57'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.64%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.06'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.38%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '143.62'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  -3.59%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.70'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.66%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.16'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.33%  '
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0483'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.56%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.11'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.81%  '
$ws.Range('E33').Value = '  -5.89%  '
$ws.Range('E34').Value = '  -2.78%  '
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('D36').Value = '1.119.24'
$ws.Range('E36').Value = '  -1.04%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.846'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -6.48%  '
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.518'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.12%  '
$ws.Range('E40').Value = '  -2.66%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '97.88'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.26%  '
$ws.Range('E42').Value = '  -4.02%  '
$ws.Range('D43').Value = '1.758.66'
$ws.Range('E43').Value = '  -1.13%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.16'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.66%  '
$ws.Range('D45').Value = '0.0₆0105'
$ws.Range('E45').Value = '  -10.10%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0530'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '54.20'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.26%  '
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.413'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('E50').Value = '  +0.37%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.46'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -4.20%  '
